# Generate Report for Handoff
# Updates status text + handoff timestamps across the three report sheets,
# and widens the "Status"/duplicate columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-21 10:45:55"
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-21 10:45:51"
$zhcn.Columns.Item(3).ColumnWidth = 16.33

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-21 10:45:55"
$dede.Columns.Item(3).ColumnWidth = 16.33
